$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.6982685
$ws.Range("H2").Value = 9.396537
$ws.Range("I2").Value = 0.1448422764790916
$ws.Range("J2").Value = 0.108319921689495
$ws.Range("Q2").Value = 1.309159988809
$ws.Range("R2").Value = 7.854959932854001
$ws.Range("S2").Value = 0.1448422764790916
$ws.Range("T2").Value = 0.108319921689495

# Row 3
$ws.Range("I3").Value = 0.1627129514638446
$ws.Range("J3").Value = 0.1825266896054461
$ws.Range("S3").Value = 0.1627129514638446
$ws.Range("T3").Value = 0.1825266896054461

# Row 4
$ws.Range("G4").Value = 5.274651666666667
$ws.Range("H4").Value = 15.823955
$ws.Range("I4").Value = 0.162611514227898
$ws.Range("J4").Value = 0.1824129002437912
$ws.Range("Q4").Value = 1.469767621178889
$ws.Range("R4").Value = 13.22790859061
$ws.Range("S4").Value = 0.162611514227898
$ws.Range("T4").Value = 0.1824129002437912

# Row 5
$ws.Range("G5").Value = 5.865133
$ws.Range("H5").Value = 11.730266
$ws.Range("I5").Value = 0.1808153824270886
$ws.Range("J5").Value = 0.1352223158932856
$ws.Range("Q5").Value = 1.634303670095334
$ws.Range("R5").Value = 9.805822020572002
$ws.Range("S5").Value = 0.1808153824270886
$ws.Range("T5").Value = 0.1352223158932856

# Row 6
$ws.Range("G6").Value = 5.310348
$ws.Range("H6").Value = 15.931044
$ws.Range("I6").Value = 0.1637119916020532
$ws.Range("J6").Value = 0.183647383979002
$ws.Range("Q6").Value = 1.479714309272
$ws.Range("R6").Value = 13.317428783448
$ws.Range("S6").Value = 0.1637119916020532
$ws.Range("T6").Value = 0.183647383979002

# Row 7
$ws.Range("G7").Value = 6.010791999999999
$ws.Range("H7").Value = 18.032376
$ws.Range("I7").Value = 0.1853058838000238
$ws.Range("J7").Value = 0.2078707885889801
$ws.Range("Q7").Value = 1.674891162021333
$ws.Range("R7").Value = 15.074020458192
$ws.Range("S7").Value = 0.1853058838000238
$ws.Range("T7").Value = 0.2078707885889801
